$p = $ppt.ActivePresentation

# 1. Update the date on the title slide (slide 1, "Date Placeholder 3").
$dateShape = $p.Slides.Item(1).Shapes.Item(3)
$dateShape.TextFrame.TextRange.Text = "2022-03-23"

# 2. Add a "Content Placeholder 2" shape with a hyperlinked test-location URL
#    to each of the four demonstration slides.
function Add-TestLocationLink($SlideIndex, $Url) {
    $slide = $p.Slides.Item($SlideIndex)

    # Re-applying the layout that is already assigned to the slide causes
    # PowerPoint to instantiate any placeholders from the layout that are
    # not yet present on the slide (here, the "Content Placeholder 2" /
    # idx=1 placeholder), without disturbing the existing Title shape.
    $slide.CustomLayout = $slide.CustomLayout

    $shape = $slide.Shapes.Item($slide.Shapes.Count)
    $textRange = $shape.TextFrame.TextRange
    $textRange.Text = $Url
    $textRange.IndentLevel = 0
    $textRange.ParagraphFormat.Bullet.Visible = 0
    $textRange.ActionSettings(1).Hyperlink.Address = $Url
}

Add-TestLocationLink 10 "https://github.com/pds-data-dictionaries/ldd-disp/tree/main/test"
Add-TestLocationLink 15 "https://github.com/pds-data-dictionaries/ldd-survey/tree/main/test"
Add-TestLocationLink 19 "https://github.com/pds-data-dictionaries/ldd-spectral/tree/main/test"
Add-TestLocationLink 24 "https://github.com/pds-data-dictionaries/ldd-nucspec/tree/main/test"
